# Commit: "merge fail - copy in files"
# The only meaningful content change in this commit is the value of cell
# B2 on the "QSfHO" worksheet (Quantization Size for Health Outcomes),
# which is updated from 0 to 1. All other differences in the target
# OOXML (theme name strings, font panose data, window geometry, calcId,
# fileVersion/rupBuild, revisionPtr GUID, and the redundant cellXfs entry)
# are artifacts produced automatically by Excel when it resaves a
# workbook and are not meaningful document edits.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsQSfHO = $wb.Worksheets.Item("QSfHO")

# Update the Quantization Size value.
$wsQSfHO.Range("B2").Value = 1

# Leave the QSfHO sheet's last-used selection on B3 (as captured in the
# saved file) without disturbing which sheet tab is active (About stays
# the active/selected tab).
$wsQSfHO.Activate()
$wsQSfHO.Range("B3").Select()
$wsAbout.Activate()
